$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; existing rows 17-28 shift down to 18-29.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = 44827
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112013
$ws.Range("G17").Value = "Alcachofa"
$ws.Range("H17").Value = "Madrigal"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 12000
$ws.Range("N17").Value = '$/caja 40 unidades'
$ws.Range("O17").Value = "Provincia de Limarí"
$ws.Range("P17").Value = 300
$ws.Range("Q17").Value = 40
$ws.Range("R17").Value = "Hortaliza"
